# Update DesignatedColors and upgrade plotting to use DesignatedColors
#
# The "Colors2" sheet holds a small table (A2:D9) of plot-series names (col A)
# together with a primary color in col B and, for rows that already had one,
# a secondary/tertiary color in col C/D.
#
# This edit designates new colors for the first three rows: the color that
# used to live in column B is shifted over into column C (making room), and a
# freshly chosen color is written into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 first, then row 2, then row 4 -- matches the order the new palette
# strings were originally added to the shared-strings table.
$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("B3").Value = "#9acd32"

$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = "#ffd700"

$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("B4").Value = "#f08080"

# Leave the selection on B9, matching the saved workbook view state.
$ws.Range("B9").Select() | Out-Null
